# Updated symbol list on Tue Dec 20 23:44:09 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple price (column D) refreshes for rows whose Coin/Link/Volume stay the same.
$priceUpdates = @{
    "D2"  = "251.73"
    "D4"  = "5.424"
    "D5"  = "0.05711"
    "D6"  = "3.421"
    "D7"  = "6.363"
    "D8"  = "0.8125"
    "D9"  = "0.9435"
    "D10" = "0.1442"
    "D11" = "0.07491"
    "D12" = "0.03170"
    "D13" = "0.03080"
    "D14" = "0.09375"
    "D15" = "3.727"
    "D16" = "0.001581"
    "D17" = "0.04764"
    "D18" = "0.0005788"
    "D19" = "0.006419"
    "D20" = "0.005043"
    "D21" = "0.001026"
    "D22" = "0.0001500"
    "D23" = "3.706"
    "D24" = "2.181"
    "D25" = "0.3304"
    "D26" = "0.1308"
    "D28" = "0.0002999"
    "D40" = "0.04030"
    "D44" = "0.008077"
    "D45" = "0.00005758"
    "D46" = "0.00000000750"
    "D47" = "0.4998"
    "D49" = "0.00002100"
    "D50" = "0.01010"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = $priceUpdates[$addr]
}

# Row 47's volume label gains a "Worstin24h" suffix this run.
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# Rows 41-43 reshuffle which coin occupies which rank (coin list order changed),
# each keeping its own B/C/D/E values as a full row record.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006770"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1072"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002710"
$ws.Range("E43").Value = "42CEJICEJI"
